$wb = $excel.ActiveWorkbook

# "Weekly Quantity" sheet: remove the week-7 order row (45361.99999999999 / 300),
# shifting all subsequent rows up by one.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(7).Delete()

# "Monthly Trend" sheet: correct the requested quantity for the
# 45382.99999999999 month row from 580 to 280.
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B4").Value = 280
